$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the cryptos.xlsx data refresh (prices / 1h volume %
# and a few re-ordered rows where coin rank changed places).
$updates = @(
    @{ Cell = "D2"; Value = "69.858.01" },
    @{ Cell = "E2"; Value = "  +2.30%  " },
    @{ Cell = "D3"; Value = "3.749.94" },
    @{ Cell = "E3"; Value = "  +19.63%  " },
    @{ Cell = "E4"; Value = "  -0.05%  " },
    @{ Cell = "D5"; Value = "617.79" },
    @{ Cell = "E5"; Value = "  +7.07%  " },
    @{ Cell = "D6"; Value = "177.82" },
    @{ Cell = "E6"; Value = "  -1.61%  " },
    @{ Cell = "D7"; Value = "3.752.23" },
    @{ Cell = "E7"; Value = "  +19.78%  " },
    @{ Cell = "E8"; Value = "  -0.04%  " },
    @{ Cell = "E9"; Value = "  +4.87%  " },
    @{ Cell = "D10"; Value = "0.169" },
    @{ Cell = "E10"; Value = "  +10.71%  " },
    @{ Cell = "D11"; Value = "6.41" },
    @{ Cell = "E11"; Value = "  -2.04%  " },
    @{ Cell = "E12"; Value = "  +7.02%  " },
    @{ Cell = "D13"; Value = "41.10" },
    @{ Cell = "E13"; Value = "  +11.80%  " },
    @{ Cell = "E14"; Value = "  +6.29%  " },
    @{ Cell = "D15"; Value = "4.374.25" },
    @{ Cell = "E15"; Value = "  +19.49%  " },
    @{ Cell = "D16"; Value = "3.750.08" },
    @{ Cell = "E16"; Value = "  +19.53%  " },
    @{ Cell = "D17"; Value = "69.987.32" },
    @{ Cell = "E17"; Value = "  +2.50%  " },
    @{ Cell = "E18"; Value = "  +1.17%  " },
    @{ Cell = "E19"; Value = "  +7.08%  " },
    @{ Cell = "D20"; Value = "518.35" },
    @{ Cell = "E20"; Value = "  +6.24%  " },
    @{ Cell = "D21"; Value = "16.80" },
    @{ Cell = "E21"; Value = "  +0.92%  " },
    @{ Cell = "D22"; Value = "9.39" },
    @{ Cell = "E22"; Value = "  +20.54%  " },
    @{ Cell = "E23"; Value = "  +6.14%  " },
    @{ Cell = "D24"; Value = "88.92" },
    @{ Cell = "E24"; Value = "  +5.98%  " },
    @{ Cell = "E25"; Value = "  +7.36%  " },
    @{ Cell = "D26"; Value = "13.59" },
    @{ Cell = "E26"; Value = "  +4.72%  " },
    @{ Cell = "D27"; Value = "10.96" },
    @{ Cell = "E27"; Value = "  +3.43%  " },
    @{ Cell = "B28"; Value = "PEPE" },
    @{ Cell = "C28"; Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe" },
    @{ Cell = "D28"; Value = "0.0000127" },
    @{ Cell = "E28"; Value = "  +33.68%  " },
    @{ Cell = "B29"; Value = "Dai" },
    @{ Cell = "C29"; Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai" },
    @{ Cell = "D29"; Value = "0.998" },
    @{ Cell = "E29"; Value = "  -0.18%  " },
    @{ Cell = "E30"; Value = "  +6.31%  " },
    @{ Cell = "E31"; Value = "  +9.02%  " },
    @{ Cell = "D32"; Value = "7.87" },
    @{ Cell = "E32"; Value = "  -3.30%  " },
    @{ Cell = "D33"; Value = "31.65" },
    @{ Cell = "E33"; Value = "  +11.83%  " },
    @{ Cell = "E34"; Value = "  +2.46%  " },
    @{ Cell = "E35"; Value = "  -0.11%  " },
    @{ Cell = "E36"; Value = "  +10.07%  " },
    @{ Cell = "D38"; Value = "0.342" },
    @{ Cell = "E38"; Value = "  +5.31%  " },
    @{ Cell = "D39"; Value = "2.20" },
    @{ Cell = "E39"; Value = "  +7.07%  " },
    @{ Cell = "E40"; Value = "  +6.71%  " },
    @{ Cell = "D41"; Value = "51.49" },
    @{ Cell = "E41"; Value = "  +4.82%  " },
    @{ Cell = "D42"; Value = "44.71" },
    @{ Cell = "E42"; Value = "  -8.72%  " },
    @{ Cell = "B43"; Value = "Cosmos" },
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom" },
    @{ Cell = "D43"; Value = "8.87" },
    @{ Cell = "E43"; Value = "  +5.17%  " },
    @{ Cell = "B44"; Value = "Bittensor" },
    @{ Cell = "C44"; Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao" },
    @{ Cell = "D44"; Value = "428.07" },
    @{ Cell = "E44"; Value = "  +8.08%  " },
    @{ Cell = "D45"; Value = "3.076.68" },
    @{ Cell = "E45"; Value = "  +10.42%  " },
    @{ Cell = "D46"; Value = "2.74" },
    @{ Cell = "E46"; Value = "  +1.31%  " },
    @{ Cell = "D47"; Value = "0.0367" },
    @{ Cell = "E47"; Value = "  +5.31%  " },
    @{ Cell = "D48"; Value = "27.91" },
    @{ Cell = "E48"; Value = "  +2.89%  " },
    @{ Cell = "B49"; Value = "ThetaToken" },
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta" },
    @{ Cell = "D49"; Value = "2.52" },
    @{ Cell = "E49"; Value = "  +7.43%  " },
    @{ Cell = "B50"; Value = "Monero" },
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr" },
    @{ Cell = "D50"; Value = "136.41" },
    @{ Cell = "E50"; Value = "  +0.57%  " },
    @{ Cell = "E51"; Value = "  -0.04%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Cell.StartsWith("D")) {
        # Column D holds price text such as "69.858.01" or "617.79" that must
        # stay a text string (it is not valid as a genuine Excel number because
        # of the thousands separators using '.'), so force text format first
        # and restore the original (unstyled) cell style afterwards.
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = $origStyle
    } else {
        $cell.Value = $u.Value
    }
}
